# Add a new "pt" (Portuguese) user row to the Users sheet and the
# corresponding computed row on the Attributes sheet.

$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users")
$wsAttrs = $wb.Worksheets.Item("Attributes")

# Users!A11:B11 - new user record
$wsUsers.Range("A11").Value = "ptuser"
$wsUsers.Range("B11").Value = "Português User"

# Attributes!A11:C11 - mirrors Users!A11 via formula, plus Group/pt attribute
$wsAttrs.Range("A11").Formula = "=Users!A11"
$wsAttrs.Range("B11").Value = "Group"
$wsAttrs.Range("C11").Value = "pt"

# Column B on Users was resized (best fit) to accommodate the new value
$wsUsers.Columns.Item(2).AutoFit()

# Leave selection/active sheet matching the state captured when the file
# was last saved: Attributes tab active, with the cell just below the new
# data selected on each sheet.
$wsUsers.Range("B12").Select()
$wsAttrs.Activate()
$wsAttrs.Range("C12").Select()
